{"js": "// Fill in the first fully-blank timesheet row (right after the\n// \"6/08/14 | 14:00-18:00 | Developed fixes for scaling and GAP requests\"\n// row) with the new entry described in the commit:\n//   Date:   6/10/14\n//   Hours:  15:30 \u2013 20:30\n//   Change: Adding support for older android versions such as kindle\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every cell's text value for every row in one batch so we can find\n// the next open (fully empty) row in the timesheet table.\nconst rowItems = rows.items;\nfor (const row of rowItems) {\n  row.cells.load(\"items/value\");\n}\nawait context.sync();\n\nlet targetRow = null;\nfor (const row of rowItems) {\n  const cellItems = row.cells.items;\n  const isFullyEmpty =\n    cellItems.length === 3 &&\n    cellItems.every((cell) => (cell.value || \"\").trim() === \"\");\n  if (isFullyEmpty) {\n    targetRow = row;\n    break;\n  }\n}\n\nif (!targetRow) {\n  throw new Error(\"Could not find an empty timesheet row to fill in.\");\n}\n\nconst [dateCell, hoursCell, changeCell] = targetRow.cells.items;\n\ndateCell.value = \"6/10/14\";\nhoursCell.value = \"15:30 \\u2013 20:30\";\nchangeCell.value = \"Adding support for older android versions such as kindle\";\n\nawait context.sync();\n", "ps1": "# Fill in the first fully-blank timesheet row (right after the\n# \"6/08/14 | 14:00-18:00 | Developed fixes for scaling and GAP requests\"\n# row) with the new entry described in the commit:\n#   Date:   6/10/14\n#   Hours:  15:30 \u2013 20:30\n#   Change: Adding support for older android versions such as kindle\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$targetRow = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $c1 = $t.Cell($r, 1).Range.Text -replace '[\\r\\a]', ''\n    $c2 = $t.Cell($r, 2).Range.Text -replace '[\\r\\a]', ''\n    $c3 = $t.Cell($r, 3).Range.Text -replace '[\\r\\a]', ''\n    if ($c1 -eq '' -and $c2 -eq '' -and $c3 -eq '') {\n        $targetRow = $r\n        break\n    }\n}\n\nif ($targetRow -eq 0) {\n    throw \"Could not find an empty timesheet row to fill in.\"\n}\n\n$t.Cell($targetRow, 1).Range.Text = \"6/10/14\"\n$t.Cell($targetRow, 2).Range.Text = \"15:30 \" + [char]0x2013 + \" 20:30\"\n$t.Cell($targetRow, 3).Range.Text = \"Adding support for older android versions such as kindle\"\n"}
